# Insert two new data rows (82 and 83) into the "Betarraga" price sheet.
# All subsequent rows (old 82..186) shift down by two rows to become 84..188.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 82:186 down by inserting two blank rows at 82:83.
$ws.Rows("82:83").Insert()

# New row 82 values.
$ws.Range("A82").Value = 9
$ws.Range("B82").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C82").Value = "Metropolitana"
$ws.Range("D82").Value = 44413
$ws.Range("E82").Value = 13
$ws.Range("F82").Value = 100114014
$ws.Range("G82").Value = "Betarraga"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 8800
$ws.Range("K82").Value = 90
$ws.Range("L82").Value = 100
$ws.Range("M82").Value = 95
$ws.Range("N82").Value = "`$/unidad"
$ws.Range("O82").Value = "Región Metropolitana"
$ws.Range("P82").Value = 95
$ws.Range("Q82").Value = 1
$ws.Range("R82").Value = "Hortaliza"

# New row 83 values.
$ws.Range("A83").Value = 9
$ws.Range("B83").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C83").Value = "Metropolitana"
$ws.Range("D83").Value = 44413
$ws.Range("E83").Value = 13
$ws.Range("F83").Value = 100114014
$ws.Range("G83").Value = "Betarraga"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Segunda"
$ws.Range("J83").Value = 3400
$ws.Range("K83").Value = 70
$ws.Range("L83").Value = 70
$ws.Range("M83").Value = 70
$ws.Range("N83").Value = "`$/unidad"
$ws.Range("O83").Value = "Región Metropolitana"
$ws.Range("P83").Value = 70
$ws.Range("Q83").Value = 1
$ws.Range("R83").Value = "Hortaliza"
